# Actualización y modificación de 3 fracciones del 4to trimestre 2020
#
# The reporting period end date ("Fecha de término del periodo que se
# informa", column C) for the two data rows was corrected from
# 2020-12-21 (serial 44186) to 2020-12-31 (serial 44196).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")
$ws.Activate()

$ws.Range("C8").Value = 44196
$ws.Range("C9").Value = 44196

# Reflect that the author had scrolled down and was working around the
# last data row when the file was saved.
$ws.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A9").Select()
